$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date value (days since 1899-12-30) that was
# bumped by one day for every data row (rows 2 through 220).
$ws.Range("C2:C220").Value = 46062
